$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "66.486.39"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -0.97%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.829.44"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +2.81%  "
$ws.Range("E4").Value = "  -0.20%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "423.07"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.46%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "130.74"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.71%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "3.822.05"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +2.82%  "
$ws.Range("E8").Value = "  -5.13%  "
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("E10").Value = "  -5.79%  "
$ws.Range("E11").Value = "  -9.31%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.0000367"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -11.25%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "40.89"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -5.06%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "4.427.59"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +2.96%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "10.12"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -5.15%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "15.62"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +16.34%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "3.833.57"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +2.88%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.138"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -1.36%  "
$ws.Range("E19").Value = "  -5.79%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "66.791.19"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -0.54%  "
$ws.Range("E21").Value = "  -6.26%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "410.96"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -8.07%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "14.48"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -11.79%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "85.37"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -5.16%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "3.03"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -4.17%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "36.87"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -2.56%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "5.67"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +12.30%  "
$ws.Range("E28").Value = "  -3.21%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "9.51"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -7.00%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "690.36"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +5.80%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "12.48"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -2.12%  "
$ws.Range("E32").Value = "  -2.61%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "2.67"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -3.57%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "7.18"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -1.79%  "
$ws.Range("E35").Value = "  -8.27%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "38.62"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -8.27%  "
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("E38").Value = "  +4.81%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "55.03"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -4.08%  "
$ws.Range("E40").Value = "  +2.53%  "
$ws.Range("E41").Value = "  -8.38%  "
$ws.Range("E42").Value = "  +0.37%  "
$ws.Range("E43").Value = "  -8.45%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "148.31"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.49%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "4.49"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +2.53%  "
$ws.Range("B46").Value = "LidoDAOToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "3.29"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -4.75%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "3.13"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -2.71%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "2.07"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -2.59%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "26.42"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -11.18%  "
$ws.Range("E50").Value = "  -4.32%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "2.55"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -4.97%  "

Write-Host "Applied cryptos update"
